# Add a new "Canada" block (NA region) to the Data sheet, mirroring the
# existing "United States of America" (NA region) rows 91:96, which share
# identical Operating Group / Sector / Vehicles Count / Miles / Accident
# Count / APMM / %Vehicles in Accidents / #Accidents with Injuries / IPMM
# values with the new Canada rows. Copy/paste preserves each cell's
# original value "kind" (text vs number) exactly - e.g. the numeric-looking
# text cells such as "0.0" stay shared-string text instead of being
# re-inferred as numbers - matching the authored diff byte-for-byte.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template rows 91:96 ("NA" / "United States of America") already contain,
# column-for-column, the same Operating Group/Sector/metrics values that
# the new Canada rows need - only the Country (column B) differs.
[void]$ws.Range("A91:K96").Copy()
[void]$ws.Range("A97:K102").PasteSpecial()

# Re-point the Country column of the newly pasted block to "Canada" (added
# to the shared-string table as a brand-new entry).
$ws.Range("B97:B102").Value = "Canada"

# Match the author's final selection recorded in the sheet view.
[void]$ws.Range("I106").Select()
